$d = $word.ActiveDocument

# The title currently reads "...Templatess": an extra duplicate "s" run
# was inserted right after "...Template" (followed by the original
# trailing "s" run). Locate the duplicated text and trim the stray
# character so the title reads "...Templates" again.
$target = $d.Content
$found = $target.Find.Execute("Templatess", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Collapse to the end of the match, then extend the start back by one
    # character so the range covers just the extra duplicated "s".
    $target.Collapse(0) | Out-Null
    $target.MoveStart(1, -1) | Out-Null
    $target.Delete() | Out-Null
}
